# Update "想去人数" (want-to-go count) figures and one cover image URL
# across the workbook's four sheets, per the upstream data refresh
# (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---- 展览 (Exhibitions) ----
$ws1.Range("F2").Value = 645
$ws1.Range("F3").Value = 722
$ws1.Range("F4").Value = 947
$ws1.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202403/OsH7V1021709288289258.jpeg"
$ws1.Range("F5").Value = 731
$ws1.Range("F6").Value = 840
$ws1.Range("F8").Value = 611
$ws1.Range("F9").Value = 134
$ws1.Range("F10").Value = 1218
$ws1.Range("F11").Value = 645
$ws1.Range("F12").Value = 388
$ws1.Range("F16").Value = 554
$ws1.Range("F18").Value = 360
$ws1.Range("F20").Value = 83
$ws1.Range("F21").Value = 557
$ws1.Range("F22").Value = 81
$ws1.Range("F23").Value = 584
$ws1.Range("F24").Value = 28
$ws1.Range("F25").Value = 794
$ws1.Range("F26").Value = 7

# ---- 演出 (Performances) ----
$ws2.Range("F4").Value = 322
$ws2.Range("F5").Value = 104
$ws2.Range("F8").Value = 184
$ws2.Range("F9").Value = 224
$ws2.Range("F10").Value = 49
$ws2.Range("F11").Value = 26
$ws2.Range("F13").Value = 99

# ---- 本地生活 (Local life) ----
$ws3.Range("F2").Value = 368

# ---- 全部类型 (All types, a merged view of the sheets above) ----
$ws4.Range("F2").Value = 368
$ws4.Range("F4").Value = 645
$ws4.Range("F6").Value = 322
$ws4.Range("F7").Value = 722
$ws4.Range("F8").Value = 947
$ws4.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202403/OsH7V1021709288289258.jpeg"
$ws4.Range("F9").Value = 731
$ws4.Range("F10").Value = 840
$ws4.Range("F12").Value = 611
$ws4.Range("F13").Value = 134
$ws4.Range("F14").Value = 1218
$ws4.Range("F15").Value = 645
$ws4.Range("F16").Value = 104
$ws4.Range("F18").Value = 388
$ws4.Range("F23").Value = 554
$ws4.Range("F24").Value = 184
$ws4.Range("F26").Value = 360
$ws4.Range("F28").Value = 83
$ws4.Range("F29").Value = 224
$ws4.Range("F30").Value = 49
$ws4.Range("F31").Value = 557
$ws4.Range("F32").Value = 26
$ws4.Range("F34").Value = 99
$ws4.Range("F35").Value = 99
$ws4.Range("F36").Value = 81
$ws4.Range("F37").Value = 584
$ws4.Range("F38").Value = 28
$ws4.Range("F39").Value = 794
$ws4.Range("F40").Value = 7
